$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# D-column values are plain numeric-looking strings (e.g. "576.88"); Excel would
# otherwise auto-convert them to numbers, so force the cell format to Text first
# to preserve them as strings, matching the original inline-string cell type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.779.69"
$ws.Range("E2").Value = "  -3.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.486.77"
$ws.Range("E3").Value = "  -5.91%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.88"
$ws.Range("E5").Value = "  -1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.25"
$ws.Range("E6").Value = "  -5.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.604"
$ws.Range("E7").Value = "  -1.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.480.93"
$ws.Range("E8").Value = "  -5.78%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  -7.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.56"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.580"
$ws.Range("E12").Value = "  -5.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.04"
$ws.Range("E13").Value = "  -6.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000270"
$ws.Range("E14").Value = "  -5.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.065.29"
$ws.Range("E15").Value = "  -5.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("E16").Value = "  -6.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "624.92"
$ws.Range("E17").Value = "  -8.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.502.41"
$ws.Range("E18").Value = "  -5.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.921.83"
$ws.Range("E19").Value = "  -3.77%  "

$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  -4.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.00"
$ws.Range("E22").Value = "  -5.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.880"
$ws.Range("E23").Value = "  -6.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.85"
$ws.Range("E24").Value = "  -9.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.47"
$ws.Range("E25").Value = "  -5.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.80"
$ws.Range("E26").Value = "  -5.05%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.61"
$ws.Range("E28").Value = "  -8.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("E29").Value = "  -11.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.38"
$ws.Range("E30").Value = "  -8.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.14"
$ws.Range("E31").Value = "  -8.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.45"
$ws.Range("E32").Value = "  -8.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  -9.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.99"
$ws.Range("E34").Value = "  -5.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "626.48"
$ws.Range("E35").Value = "  +6.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.69"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.47"
$ws.Range("E37").Value = "  -15.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -6.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0443"
$ws.Range("E41").Value = "  -3.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("E42").Value = "  -6.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.349.02"
$ws.Range("E43").Value = "  -9.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("E44").Value = "  -6.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.61"
$ws.Range("E45").Value = "  -8.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0687"
$ws.Range("E46").Value = "  -11.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").Value = "  -8.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  -3.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.28"
$ws.Range("E50").Value = "  -3.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.57"
$ws.Range("E51").Value = "  +13.02%  "

# Rows 39 and 40: OKB and FirstDigitalUSD swap rank positions
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.52"
$ws.Range("E40").Value = "  -4.53%  "
